{"js": "// Add a left indent of 567 twips (28.35 pt / 1 cm) with a hanging indent\n// of 283 twips (14.15 pt / 0.5 cm) to every \"List Paragraph\" (style \"ae\")\n// paragraph in the \"\u0388\u03c7\u03bf\u03bd\u03c4\u03b1\u03c2 \u03c5\u03c0\u03cc\u03c8\u03b7\" (Having regard to) numbered list.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.style === \"List Paragraph\") {\n    paragraph.leftIndent = 28.35; // 567 twips\n    paragraph.firstLineIndent = -14.15; // -283 twips (hanging indent)\n  }\n}\n\nawait context.sync();\n", "ps1": "# Add a left indent of 567 twips (28.35 pt / 1 cm) with a hanging indent\n# of 283 twips (14.15 pt / 0.5 cm) to every \"List Paragraph\" (style \"ae\")\n# paragraph in the \"\u0388\u03c7\u03bf\u03bd\u03c4\u03b1\u03c2 \u03c5\u03c0\u03cc\u03c8\u03b7\" (Having regard to) numbered list.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"List Paragraph\") {\n        $p.Format.LeftIndent = 28.35\n        $p.Format.FirstLineIndent = -14.15\n    }\n}\n"}
